# Updated cryptos list on Wed Jan 24 11:55:29 UTC 2024 with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) values for
# each coin row on the active worksheet with the latest scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.176.94"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.243.01"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.61%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.02"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.67%  "

$ws.Range("E7").Value = "  +2.72%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "31.19"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +12.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.27"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.91%  "

$ws.Range("E14").Value = "  +6.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.591.52"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.253.94"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.739"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "40.081.64"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.69%  "

$ws.Range("E20").Value = "  +4.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.71"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +9.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.76"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.56"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.74%  "

$ws.Range("E25").Value = "  +0.13%  "

$ws.Range("E26").Value = "  +3.90%  "

$ws.Range("E27").Value = "  +8.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.03"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.64%  "

$ws.Range("E30").Value = "  +5.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.94"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.81%  "

$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("E34").Value = "  +3.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0720"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.84%  "

$ws.Range("E36").Value = "  +3.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.63"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +16.45%  "

$ws.Range("E38").Value = "  +6.94%  "

$ws.Range("E39").Value = "  +3.57%  "

$ws.Range("E40").Value = "  +3.66%  "

$ws.Range("E41").Value = "  +7.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.83"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.027.48"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.21"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +10.99%  "

$ws.Range("E45").Value = "  +7.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.98"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +11.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("E48").Value = "  +3.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.471.87"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.66"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.13"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.11%  "
